# Add data for 2022-12-04: bump the "through" date from 11-25 to 11-26,
# and update the new November / Total figures for the 2022 column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-11-26"

# Update the 2022 column header label.
$ws.Range("I1").Value = "2022 (through 11-26)"

# Update November (row 12) and Total (row 14) figures for the 2022 column.
$ws.Range("I12").Value = 97
$ws.Range("I14").Value = 1495
